# Figure 5 "Data access" diagram (slide 1): reword 'variable' -> 'descriptor'
# in two process-flow shapes, matching the upstream commit
# "Changed 'variable' to 'descriptor' in Fig 5".
#
#   - "Rounded Rectangle 51": single-run text
#         "Identify variable(s) of interest"
#     becomes
#         "Identify descriptor(s) of interest"
#
#   - "Rounded Rectangle 52": only the FIRST run changes
#         "Download variable(s) from "  ->  "Download descriptor(s) from "
#     the remaining runs ("Zenodo", " and decompress tiles") and their
#     run-level formatting must stay untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp1 = $s.Shapes.Item("Rounded Rectangle 51")
$shp1.TextFrame.TextRange.Runs(1).Text = "Identify descriptor(s) of interest"

$shp2 = $s.Shapes.Item("Rounded Rectangle 52")
$shp2.TextFrame.TextRange.Runs(1).Text = "Download descriptor(s) from "
